# "Sadia sister is not study" -> remove her from the routine.
# The cells I6, I8 and I10 held the shared string "Sadia Sister";
# clearing them removes the last references to that string so it
# drops out of the shared-strings table entirely on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I6").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("I10").Value = ""

$null = $ws.Range("I10").Select()
